$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("Golang Developer", "https://www.dice.com/job-detail/a07931e6-29b6-4e82-b0c7-9ce1edd1fde7", "Fremont, California", "Full-time, Contract", "Depends on Experience", "Radyant Inc."),
    @("Software Engineer - GO/JavaScript/Angular Or React", "https://www.dice.com/job-detail/de55429f-3edf-46c3-a6b6-05601d6ccbd4", "Remote", "Contract", "Depends on Experience", "Ocean Blue Solution"),
    @("Golang Go Architect", "https://www.dice.com/job-detail/7afb0fd0-b3ca-4475-a396-1c2a8b13513a", "Atlanta, Georgia", "Contract", "$80 - $85", "Source Mantra Inc"),
    @("Golang Developer", "https://www.dice.com/job-detail/f56507c9-5983-470e-a67a-e683013b4296", "Remote", "Contract", "$90 - $100", "Arnex Solutions LLC")
)

$startRow = 39
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowNum = $startRow + $i
    $rowData = $newRows[$i]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $rowData[$col - 1]
    }
}
